# Automatische test-sync: 2025-08-14 22:06:50
# Append the new "Klacht over levering" log entry to the Logs sheet (row 38),
# extend the conditional-formatting ranges that covered the previous last
# row (37) so they include the new row 38, and bump the Dashboard summary
# count for "Intern verzoek / Actie voor medewerker" from 29 to 30.

$wb = $excel.ActiveWorkbook

# --- 1. Append new row to the "Logs" sheet -------------------------------
$logs = $wb.Worksheets.Item("Logs")

$newRow = 38
$logs.Cells.Item($newRow, 1).Value  = "Klacht over levering"
$logs.Cells.Item($newRow, 2).Value  = "klantenservice@testbedrijf123.nl"
$logs.Cells.Item($newRow, 3).Value  = "Ik ben niet tevreden over hoe dit is gegaan."
$logs.Cells.Item($newRow, 4).Value  = "Intern verzoek / Actie voor medewerker"
$logs.Cells.Item($newRow, 5).Value  = "Bedankt, we hebben dit doorgestuurd naar klachten@testbedrijf123.nl."
$logs.Cells.Item($newRow, 6).Value  = "2025-08-14 22:06:29"
$logs.Cells.Item($newRow, 7).Value  = "Nee"
$logs.Cells.Item($newRow, 8).Value  = "Ja"
$logs.Cells.Item($newRow, 9).Value  = "Nee"
$logs.Cells.Item($newRow, 10).Value = "Nee"

# --- 2. Extend conditional formatting ranges from row 37 to row 38 -------
$cfColumns = @("D", "G", "H", "I", "J")
foreach ($col in $cfColumns) {
    $oldRange = $logs.Range($col + "2:" + $col + "37")
    $newRange = $logs.Range($col + "2:" + $col + "38")
    $fcs = $oldRange.FormatConditions
    $count = $fcs.Count
    for ($i = 1; $i -le $count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- 3. Update the Dashboard summary count --------------------------------
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 30

Write-Output "Logs row 38 added, conditional formatting extended, Dashboard B2 updated."
